$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Adults" column (E) values were bumped by 1 and are now stored as text
# instead of numbers (rows 2-5: 1 -> "2", rows 6-9: 2 -> "3").
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 5).Value = "2"
}
for ($r = 6; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = "3"
}

# Selection moved from A1 to E10.
$ws.Range("E10").Select()
